$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values are kept as text (matching the original
# inlineStr cell type) so numeric-looking strings (e.g. "26.00", "0.0000118")
# aren't silently coerced/re-formatted as numbers.
$updates = [ordered]@{
    'D2' = '64.295.84'
    'E2' = '  -2.43%  '
    'D3' = '3.179.50'
    'E3' = '  -7.62%  '
    'E4' = '  +0.00%  '
    'D5' = '563.30'
    'D6' = '171.68'
    'E6' = '  -1.19%  '
    'E7' = '  +0.03%  '
    'D8' = '0.609'
    'E8' = '  +1.68%  '
    'D9' = '3.176.36'
    'E9' = '  -7.64%  '
    'E10' = '  -5.53%  '
    'D11' = '6.62'
    'E11' = '  -4.39%  '
    'E12' = '  -3.25%  '
    'D13' = '3.734.22'
    'E13' = '  -7.51%  '
    'E14' = '  +1.38%  '
    'D15' = '27.52'
    'E15' = '  -4.74%  '
    'D16' = '64.290.62'
    'E16' = '  -2.42%  '
    'E17' = '  -4.68%  '
    'D18' = '3.177.91'
    'E18' = '  -7.70%  '
    'E19' = '  -4.02%  '
    'D20' = '13.11'
    'E20' = '  -4.64%  '
    'D21' = '354.03'
    'E21' = '  -4.31%  '
    'D22' = '7.22'
    'E22' = '  -4.83%  '
    'E23' = '  +0.39%  '
    'D24' = '69.14'
    'E24' = '  -4.22%  '
    'D25' = '0.506'
    'E25' = '  -4.20%  '
    'D26' = '0.0000118'
    'E26' = '  -2.77%  '
    'D27' = '9.72'
    'E27' = '  +0.32%  '
    'D28' = '0.175'
    'E28' = '  -1.50%  '
    'E29' = '  +0.12%  '
    'E30' = '  -1.17%  '
    'E31' = '  -0.04%  '
    'E32' = '  -3.66%  '
    'D33' = '22.13'
    'E33' = '  -6.13%  '
    'D34' = '6.67'
    'E34' = '  -4.39%  '
    'E35' = '  -5.41%  '
    'D36' = '1.45'
    'E36' = '  -5.22%  '
    'D37' = '155.37'
    'E37' = '  -3.46%  '
    'D38' = '0.810'
    'E38' = '  -7.81%  '
    'D39' = '26.00'
    'E39' = '  -8.10%  '
    'E40' = '  -4.23%  '
    'E41' = '  -4.73%  '
    'D42' = '2.653.60'
    'E42' = '  -4.08%  '
    'E43' = '  -6.11%  '
    'B44' = 'Bittensor'
    'C44' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D44' = '334.10'
    'E44' = '  +3.39%  '
    'B45' = 'RenderToken'
    'C45' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D45' = '6.02'
    'E45' = '  -6.86%  '
    'E46' = '  -3.54%  '
    'D47' = '38.61'
    'E47' = '  -3.54%  '
    'D48' = '23.79'
    'E48' = '  -3.65%  '
    'E49' = '  -5.96%  '
    'E50' = '  -0.93%  '
    'E51' = '  -0.04%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = '@'
    $range.Value = $updates[$cellRef]
}
